$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$active = $wb.ActiveSheet
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet4"
$active.Activate()
